# Reorder the metric_inputs_adults workbook so that, per the Generalized
# Network Code convention, the "Path Transition" block comes before the
# "Path Survival" block.  Rows 6-10 and rows 11-15 (labels, values,
# formulas and the small number inputs that go with them) trade places,
# and the blank "spacer" row above each block (row 5 / row 10) swaps its
# distinguishing alignment formatting along with its block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the two content blocks as they exist today ----------------
# Block A: rows 6-9   ("Path Survival" section, old layout)
$A_B6  = $ws.Range("B6").Value()
$A_A7  = $ws.Range("A7").Value()
$A_C7  = $ws.Range("C7").Value()
$A_D7  = $ws.Range("D7").Value()
$A_A8  = $ws.Range("A8").Value()
$A_B8  = $ws.Range("B8").Value()
$A_C8  = $ws.Range("C8").Value()
$A_D8  = $ws.Range("D8").Value()
$A_B9  = $ws.Range("B9").Value()
$A_C9  = $ws.Range("C9").Value()
$A_D9  = $ws.Range("D9").Value()

# Block B: rows 11-14 ("Path Transition" section, old layout)
$B_B11 = $ws.Range("B11").Value()
$B_A12 = $ws.Range("A12").Value()
$B_C12 = $ws.Range("C12").Value()
$B_D12 = $ws.Range("D12").Value()
$B_A13 = $ws.Range("A13").Value()
$B_B13 = $ws.Range("B13").Value()
$B_F13 = $ws.Range("F13").Value()
$B_B14 = $ws.Range("B14").Value()
$B_C14 = $ws.Range("C14").Value()
$B_D14 = $ws.Range("D14").Value()

# --- clear the cells whose address isn't reused by the other block -----
$ws.Range("F13").ClearContents()
$ws.Range("E14").ClearContents()

# --- Block B's content moves up into rows 6-9 (new first section) ------
$ws.Range("B6").Value = $B_B11
$ws.Range("A7").Value = $B_A12
$ws.Range("C7").Value = $B_C12
$ws.Range("D7").Value = $B_D12
$ws.Range("A8").Value = $B_A13
$ws.Range("B8").Value = $B_B13
$ws.Range("C8").Formula = "=F8/C3*(D3-1)"
$ws.Range("D8").Formula = "=1-C8"
$ws.Range("F8").Value = $B_F13
$ws.Range("B9").Value = $B_B14
$ws.Range("C9").Value = $B_C14
$ws.Range("D9").Value = $B_D14
$ws.Range("E9").Formula = "=C3*D3"

# --- Block A's content moves down into rows 11-14 (new second section) --
$ws.Range("B11").Value = $A_B6
$ws.Range("A12").Value = $A_A7
$ws.Range("C12").Value = $A_C7
$ws.Range("D12").Value = $A_D7
$ws.Range("A13").Value = $A_A8
$ws.Range("B13").Value = $A_B8
$ws.Range("C13").Value = $A_C8
$ws.Range("D13").Value = $A_D8
$ws.Range("B14").Value = $A_B9
$ws.Range("C14").Value = $A_C9
$ws.Range("D14").Value = $A_D9

# --- the trailing helper formula on row 15 referenced E14; now it must
#     reference E9, since that's where "(100 - used up in transitions)"
#     now lives ---------------------------------------------------------
$ws.Range("E15").Formula = "=(100-E9)/C4"

# --- the blank spacer rows (5 and 10) swap their alignment formatting so
#     the "new section" visual cue travels with the block that now sits
#     below it ------------------------------------------------------------
$ws.Range("ZZ1").Value = "tmp"
$ws.Range("A10").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)   # xlPasteFormats -> ZZ1 holds A10's original look

$ws.Range("A5").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # A10 <- A5's original look

$ws.Range("ZZ1").Copy()
$ws.Range("A5").PasteSpecial(-4122)    # A5 <- A10's original look (via ZZ1)

$ws.Range("ZZ1").Clear()

# --- tidy the selection / view state to match the post-edit workbook ---
$ws.Range("E23").Select()
